# Changed duration of tutorials to 1hr
# - Section_A and Section_B sheets: re-shuffled class slots and split the
#   old 1.5hr afternoon slot (12:30-17:00 leftover block) into dedicated
#   1hr slots, adding 5 new rows (8-12) to each timetable.
# - Course_Summary sheet: replaced the single "Credits" column with
#   "Lectures/Week", "Tutorials/Week" and "Total Credits" columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Section_A
# ---------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("B2").Value = "Free"
$wsA.Range("E2").Value = "CS264"
$wsA.Range("F2").Value = "Free"

$wsA.Range("B3").Value = "CS263"
$wsA.Range("C3").Value = "CS263"
$wsA.Range("D3").Value = "Free"
$wsA.Range("E3").Value = "MA261"
$wsA.Range("F3").Value = "CS264"

$wsA.Range("D5").Value = "CS263"
$wsA.Range("F5").Value = "CS261"

$wsA.Range("B6").Value = "CS264"
$wsA.Range("C6").Value = "MA261"
$wsA.Range("F6").Value = "Free"

$wsA.Range("C7").Value = "Free"

# New rows 8-12 (copy formatting of the existing time-slot row 7 first,
# so the new time-slot cells in column A pick up the same bold/border
# style used by A2:A7)
$wsA.Range("A8:F12").Value = "Free"

$wsA.Range("A8").Value = "12:00-13:00"
$wsA.Range("A9").Value = "13:00-14:00"
$wsA.Range("A10").Value = "15:30-16:30"
$wsA.Range("A11").Value = "16:30-17:30"
$wsA.Range("A12").Value = "17:30-18:30"

$wsA.Range("F12").Value = "CS264 (Tutorial)"

$wsA.Range("A7").Copy()
$wsA.Range("A8:A12").PasteSpecial(-4122)
$wsA.Application.CutCopyMode = $false

# ---------------------------------------------------------------
# Section_B
# ---------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("C2").Value = "MA261"

$wsB.Range("C3").Value = "CS261"
$wsB.Range("D3").Value = "CS261"
$wsB.Range("E3").Value = "Free"

$wsB.Range("B5").Value = "Free"
$wsB.Range("C5").Value = "CS264"
$wsB.Range("D5").Value = "Free"
$wsB.Range("E5").Value = "CS263"
$wsB.Range("F5").Value = "CS263"

$wsB.Range("B6").Value = "CS264"
$wsB.Range("C6").Value = "CS263"
$wsB.Range("D6").Value = "Free"
$wsB.Range("E6").Value = "Free"

$wsB.Range("B7").Value = "CS261"
$wsB.Range("C7").Value = "Free"
$wsB.Range("D7").Value = "MA261"
$wsB.Range("E7").Value = "Free"
$wsB.Range("F7").Value = "CS264"

# New rows 8-12
$wsB.Range("A8:F12").Value = "Free"

$wsB.Range("A8").Value = "12:00-13:00"
$wsB.Range("A9").Value = "13:00-14:00"
$wsB.Range("A10").Value = "15:30-16:30"
$wsB.Range("A11").Value = "16:30-17:30"
$wsB.Range("A12").Value = "17:30-18:30"

$wsB.Range("E9").Value = "CS264 (Tutorial)"

$wsB.Range("A7").Copy()
$wsB.Range("A8:A12").PasteSpecial(-4122)
$wsB.Application.CutCopyMode = $false

# ---------------------------------------------------------------
# Course_Summary
# ---------------------------------------------------------------
$wsC = $wb.Worksheets.Item("Course_Summary")

# Split the old "Credits" column (F) into three columns: insert two new
# columns before the Instructor column so it becomes H, keeping F & G
# for the new "Lectures/Week" / "Tutorials/Week" values and adding a
# "Total Credits" header in G (old Instructor column shifts right).
$wsC.Range("F1:G1").EntireColumn.Insert()

$wsC.Range("E1").Value = "Lectures/Week"
$wsC.Range("F1").Value = "Tutorials/Week"
$wsC.Range("G1").Value = "Total Credits"
$wsC.Range("H1").Value = "Instructor"

# MA261
$wsC.Range("E2").Value = 2
$wsC.Range("F2").Value = 0
$wsC.Range("G2").Value = 2

# CS261
$wsC.Range("E3").Value = 3
$wsC.Range("F3").Value = 0
$wsC.Range("G3").Value = 2

# CS263
$wsC.Range("E4").Value = 3
$wsC.Range("F4").Value = 0
$wsC.Range("G4").Value = 4

# CS264
$wsC.Range("E5").Value = 3
$wsC.Range("F5").Value = 1
$wsC.Range("G5").Value = 4
